$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bugfix: the access-level/user rows were out of sync with their accounts.
# Re-align rows 3-6 so each user keeps their correct access level, and
# rotate Alan Moreno to the bottom of the list with the corrected level.
$ws.Range("A3").Value = "Scott Luzader"
$ws.Range("C3").Value = 3

$ws.Range("A4").Value = "Paolo Cisneros"
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = "Samuel Lopez"
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = "Alan Moreno"
$ws.Range("C6").Value = 2

# Update the active selection to match the saved view.
$ws.Range("F3").Select()
